$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,7).Value = 22.59487733333333
$ws.Cells.Item(2,8).Value = 67.784632
$ws.Cells.Item(2,9).Value = 0.7395019553569895
$ws.Cells.Item(2,10).Value = 0.7395019553569895
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,13).Value = 34.07074633333333
$ws.Cells.Item(2,14).Value = 102.212239
$ws.Cells.Item(2,15).Value = 0.5171464495142372
$ws.Cells.Item(2,16).Value = 0.5171464495142373
$ws.Cells.Item(2,17).Value = 769.824334056783
$ws.Cells.Item(2,18).Value = 6928.419006511047
$ws.Cells.Item(2,19).Value = 0.382430810621703
$ws.Cells.Item(2,20).Value = 0.3824308106217031

# Row 3
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,7).Value = 22.59487733333333
$ws.Cells.Item(3,8).Value = 67.784632
$ws.Cells.Item(3,9).Value = 0.7395019553569895
$ws.Cells.Item(3,10).Value = 0.7395019553569895
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,13).Value = 27.685497
$ws.Cells.Item(3,14).Value = 83.056491
$ws.Cells.Item(3,15).Value = 0.420227262899125
$ws.Cells.Item(3,16).Value = 0.4202272628991251
$ws.Cells.Item(3,17).Value = 625.5504086273679
$ws.Cells.Item(3,18).Value = 5629.953677646312
$ws.Cells.Item(3,19).Value = 0.3107588826082187
$ws.Cells.Item(3,20).Value = 0.3107588826082187

# Row 4
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,7).Value = 22.59487733333333
$ws.Cells.Item(4,8).Value = 67.784632
$ws.Cells.Item(4,9).Value = 0.7395019553569895
$ws.Cells.Item(4,10).Value = 0.7395019553569895
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,13).Value = 4.125957666666666
$ws.Cells.Item(4,14).Value = 12.377873
$ws.Cells.Item(4,15).Value = 0.06262628758663766
$ws.Cells.Item(4,16).Value = 0.06262628758663766
$ws.Cells.Item(4,17).Value = 93.22550736085955
$ws.Cells.Item(4,18).Value = 839.0295662477359
$ws.Cells.Item(4,19).Value = 0.0463122621270677
$ws.Cells.Item(4,20).Value = 0.0463122621270677

# Row 5
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,7).Value = 0.3045986666666667
$ws.Cells.Item(5,8).Value = 0.913796
$ws.Cells.Item(5,9).Value = 0.00996913177602551
$ws.Cells.Item(5,10).Value = 0.00996913177602551
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,13).Value = 34.07074633333333
$ws.Cells.Item(5,14).Value = 102.212239
$ws.Cells.Item(5,15).Value = 0.5171464495142372
$ws.Cells.Item(5,16).Value = 0.5171464495142373
$ws.Cells.Item(5,17).Value = 10.37790390547155
$ws.Cells.Item(5,18).Value = 93.40113514924398
$ws.Cells.Item(5,19).Value = 0.005155501102711154
$ws.Cells.Item(5,20).Value = 0.005155501102711155

# Row 6
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,7).Value = 0.3045986666666667
$ws.Cells.Item(6,8).Value = 0.913796
$ws.Cells.Item(6,9).Value = 0.00996913177602551
$ws.Cells.Item(6,10).Value = 0.00996913177602551
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,13).Value = 27.685497
$ws.Cells.Item(6,14).Value = 83.056491
$ws.Cells.Item(6,15).Value = 0.420227262899125
$ws.Cells.Item(6,16).Value = 0.4202272628991251
$ws.Cells.Item(6,17).Value = 8.432965472204
$ws.Cells.Item(6,18).Value = 75.896689249836
$ws.Cells.Item(6,19).Value = 0.004189300959719894
$ws.Cells.Item(6,20).Value = 0.004189300959719894

# Row 7
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,7).Value = 0.3045986666666667
$ws.Cells.Item(7,8).Value = 0.913796
$ws.Cells.Item(7,9).Value = 0.00996913177602551
$ws.Cells.Item(7,10).Value = 0.00996913177602551
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,13).Value = 4.125957666666666
$ws.Cells.Item(7,14).Value = 12.377873
$ws.Cells.Item(7,15).Value = 0.06262628758663766
$ws.Cells.Item(7,16).Value = 0.06262628758663766
$ws.Cells.Item(7,17).Value = 1.256761203989778
$ws.Cells.Item(7,18).Value = 11.310850835908
$ws.Cells.Item(7,19).Value = 0.0006243297135944615
$ws.Cells.Item(7,20).Value = 0.0006243297135944615

# Row 8
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,7).Value = 7.654706
$ws.Cells.Item(8,8).Value = 22.964118
$ws.Cells.Item(8,9).Value = 0.2505289128669849
$ws.Cells.Item(8,10).Value = 0.2505289128669849
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,13).Value = 34.07074633333333
$ws.Cells.Item(8,14).Value = 102.212239
$ws.Cells.Item(8,15).Value = 0.5171464495142372
$ws.Cells.Item(8,16).Value = 0.5171464495142373
$ws.Cells.Item(8,17).Value = 260.8015463822446
$ws.Cells.Item(8,18).Value = 2347.213917440201
$ws.Cells.Item(8,19).Value = 0.1295601377898229
$ws.Cells.Item(8,20).Value = 0.129560137789823

# Row 9
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,7).Value = 7.654706
$ws.Cells.Item(9,8).Value = 22.964118
$ws.Cells.Item(9,9).Value = 0.2505289128669849
$ws.Cells.Item(9,10).Value = 0.2505289128669849
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,13).Value = 27.685497
$ws.Cells.Item(9,14).Value = 83.056491
$ws.Cells.Item(9,15).Value = 0.420227262899125
$ws.Cells.Item(9,16).Value = 0.4202272628991251
$ws.Cells.Item(9,17).Value = 211.924339998882
$ws.Cells.Item(9,18).Value = 1907.319059989938
$ws.Cells.Item(9,19).Value = 0.1052790793311865
$ws.Cells.Item(9,20).Value = 0.1052790793311865

# Row 10
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,7).Value = 7.654706
$ws.Cells.Item(10,8).Value = 22.964118
$ws.Cells.Item(10,9).Value = 0.2505289128669849
$ws.Cells.Item(10,10).Value = 0.2505289128669849
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,13).Value = 4.125957666666666
$ws.Cells.Item(10,14).Value = 12.377873
$ws.Cells.Item(10,15).Value = 0.06262628758663766
$ws.Cells.Item(10,16).Value = 0.06262628758663766
$ws.Cells.Item(10,17).Value = 31.58299290677933
$ws.Cells.Item(10,18).Value = 284.246936161014
$ws.Cells.Item(10,19).Value = 0.01568969574597549
$ws.Cells.Item(10,20).Value = 0.01568969574597549
